$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 100, shifting existing rows 100-132 down to 101-133
$ws.Rows.Item(100).Insert()

# Populate the new row 100 with the weekly update's data
$ws.Range("A100").Value = 1
$ws.Range("B100").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C100").Value = "Arica y Parinacota"
$ws.Range("D100").Value = 45215
$ws.Range("E100").Value = 15
$ws.Range("F100").Value = 100112040
$ws.Range("G100").Value = "Cilantro"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 580
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = 1241
$ws.Range("N100").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 827
$ws.Range("Q100").Value = 1.5
$ws.Range("R100").Value = "Hortaliza"

# Keep the date column formatted like the rest of column D
$ws.Range("D100").NumberFormat = $ws.Range("D101").NumberFormat
